$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Target cluster" for the (only remaining) data row changes from
# "ECs" to "FAPs", and its metrics are refreshed with the newly
# recomputed TPM-based values (previously held by the row that is being
# removed below).
$ws.Range("D2").Value = "FAPs"

$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.006582333333333333
$ws.Range("N2").Value = 0.019747
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.00009554476244444445
$ws.Range("R2").Value = 0.0008599028620000001
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Remove the old duplicate "FAPs" data row (row 3); its figures now live
# on row 2.
$ws.Rows("3:3").Delete()
